# The commit swaps the full data payload (columns B..AD) between pairs of
# adjacent rows in the "Turkey Super Lig" sheet, while leaving column A
# (the running row index) untouched. Each pair below corresponds to two
# fixtures whose rows had their B:AD contents exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(32, 33),
    @(50, 51),
    @(101, 102),
    @(112, 113),
    @(166, 167),
    @(171, 172),
    @(197, 198),
    @(199, 200),
    @(205, 206),
    @(208, 209),
    @(230, 231),
    @(241, 242),
    @(243, 244),
    @(245, 246),
    @(270, 271),
    @(323, 324),
    @(363, 364),
    @(373, 374)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # NOTE: use ${r1}/${r2} delimiters -- "B$r1:AD$r1" would be mis-parsed
    # by PowerShell as a drive-qualified variable reference ("$r1:"), which
    # silently swallows everything from the colon up to the next space.
    $range1 = $ws.Range("B${r1}:AD${r1}")
    $range2 = $ws.Range("B${r2}:AD${r2}")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
